$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.677.11"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.509.15"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "618.90"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.94"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.503.65"
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.197"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.576"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "45.39"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000272"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.079.82"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.33"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "597.35"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.518.07"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.792.09"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.50"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.872"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.98"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.46"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.95"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.68"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.20"
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.88"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.95"
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.70"
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "596.69"
$ws.Range("E35").Value = "  -16.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0496"
$ws.Range("E36").Value = "  +4.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.78"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0982"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.76"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.143"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.308.75"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0715"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.305"
$ws.Range("E45").Value = "  -4.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.83"
$ws.Range("E46").Value = "  -3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.35"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.45"
$ws.Range("E48").Value = "  -5.57%  "
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.71"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("E51").Value = "  -0.05%  "
